# Update 18-Jun-2021, end of day update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 26 (18-Jun-2021 entries continue under existing date 44364) ---
$ws.Range("B26").Value = "Wages Expense"
$ws.Range("D26").Formula = "=60000+280000"

# --- Row 27 ---
$ws.Range("B27").Value = "TRANSFER BCA"
$ws.Range("D27").Formula = "=3450000+2600000+2625000+2750000+12840000+9000000+2690000+1608000"

# --- Row 28 ---
$ws.Range("B28").Value = "A/P"
$ws.Range("C28").Formula = "=34938000"

# --- Row 29 ---
$ws.Range("B29").Value = "SALES - cash/retail"
$ws.Range("C29").Formula = "=37378275+5255725-34938000"

# --- Row 30 ---
$ws.Range("B30").Value = "SELISIH - lebih"
$ws.Range("C30").Value = 30000

# --- Row 31 ---
$ws.Range("B31").Value = "SETOR KE BANK"
$ws.Range("D31").Value = 4000000

# --- Row 32 (new date 19-Jun-2021, serial 44365) ---
$ws.Range("A32").Value = 44365
$ws.Range("A32").NumberFormat = $ws.Range("A31").NumberFormat
$ws.Range("B32").Value = "Wages Expense"
$ws.Range("D32").Formula = "=60000"

# --- Row 33 ---
$ws.Range("B33").Value = "TRANSFER BCA"
$ws.Range("D33").Formula = "=775000+577200+1555000"

# --- Row 34 ---
$ws.Range("B34").Value = "BELI lampu"
$ws.Range("D34").Value = 32500

# --- Row 35 ---
$ws.Range("B35").Value = "A/P"
$ws.Range("C35").Formula = "=577200"

# --- Row 36 ---
$ws.Range("B36").Value = "FREIGHT OUT"
$ws.Range("D36").Formula = "=162500"

# --- Update the view: scroll frozen pane and move active cell/selection ---
$ws.Application.ActiveWindow.SelectedSheets.Item(1).Select()
$ws.Range("B37").Select()
$ws.Application.ActiveWindow.ScrollRow = 18
